$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: F1_Score label and formula
$ws.Range("A10").Value = "F1_Score"
$ws.Range("B10").Formula = "=(2*B8*B9)/(B8+B9)"

# Match the percentage style used by B8/B9
$ws.Range("B10").NumberFormat = $ws.Range("B8").NumberFormat

# Update selection to match target state
$ws.Range("C12").Select()
